# TAK update, Psychiatry addition
#
# - Diseases screen (Main sheet):
#     * F5 "Genetic Diseases" is replaced by a new "Gastroenterology" entry
#     * the old "Genetic Diseases" entry is re-added one row down (F6)
#     * Psychiatry (F12) gets a TAK/last-checked date in G12
# - Largest (col I/J) TAK date (J4) is refreshed
# - Active cell/selection moves to G13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh "Largest" TAK date
$ws.Range("J4").Value2 = 45054

# Diseases list: swap Genetic Diseases -> Gastroenterology, push Genetic Diseases down a row
$ws.Range("F5").Value2 = "Gastroenterology"
$ws.Range("F6").Value2 = "Genetic Diseases"

# Psychiatry: add TAK date, matching the date formatting used by sibling cells in column G
$ws.Range("G12").Value2 = 45067
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the author left it
$ws.Range("G13").Select()
